$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new (blank) rows at position 13, pushing the existing
# "Programa resumido:" ... "Requisitos:" block down to rows 15..24
# and the two requisito rows down to 26/27. The inserted rows pick up
# the sheet's default (15pt) height, which is exactly what the target
# layout needs for every later row.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# The two new rows only carry B/C data (no label in column A), so drop
# whatever got carried into column A by the insert.
$ws.Range("A13:A14").Clear()

# Give the new B13:C14 cells the same look (wrap text / red "modified"
# font) as the rest of the B/C columns by copying the formatting from
# a neighboring data row.
$ws.Range("B15:C15").Copy()
$ws.Range("B13:C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Objetivos (row 10): replace the instructor name that had been
# pasted into B/C with the actual course-objectives text.
$ws.Range("B10").Value2 = "Ensinar a identificação e especificação dos elementos que compõem as tubulações que integram os processos inerentes às industrias de processamento.Auxiliar o desenvolvimento da habilidade de planejamento e projeto de processos industriais."
$ws.Range("C10").Value2 = $ws.Range("B10").Value2

# --- Docentes responsáveis (rows 13/14): one professor per row.
$ws.Range("B13").Value2 = "6634418 - Antonio Clelio Ribeiro"
$ws.Range("C13").Value2 = $ws.Range("B13").Value2

$ws.Range("B14").Value2 = "1285870 - Marcos Villela Barcza"
$ws.Range("C14").Value2 = $ws.Range("B14").Value2

# --- Programa resumido (row 15): real short-syllabus text instead of
# the stray activation date.
$ws.Range("B15").Value2 = "Tubos e Tubulações - DefiniçõesTubos: Materiais, Processos de Fabricação e Normalização Dimensional.Meios de Ligação de Tubos, Conexões de Tubulações e Juntas de Expansão.VálvulasPurgadores de Vapor, Separadores e Filtros. Recomendações de Material para Serviços. Aquecimento, Isolamento Térmico, Pintura e Proteção.Disposição das Construções em uma Instalação Industrial. Arranjo e Detalhamento de Tubulações.Sistemas Especiais de Tubulação. Suportes de Tubulação. Montagem e Teste de Tubulações.Visita Técnica Desenhos de TubulaçõesCálculo do diâmetro das tubulaçõesA Tubulação Considerada como Elemento Estrutural Cálculo da Espessura de Parede de Tubos e do Vão entre Suportes.Dilatação Térmica e Flexibilidade de Tubulações. Cálculo de Cálculo de Flexibilidade."
$ws.Range("C15").Value2 = $ws.Range("B15").Value2

# --- Programa (row 17): full syllabus text instead of the stray
# instructor name.
$ws.Range("B17").Value2 = "Tubos e Tubulações - DefiniçõesTubos: Materiais, Processos de Fabricação e Normalização Dimensional.Meios de Ligação de Tubos, Conexões de Tubulações e Juntas de Expansão.VálvulasPurgadores de Vapor, Separadores e Filtros. Recomendações de Material para Serviços. Aquecimento, Isolamento Térmico, Pintura e Proteção.Disposição das Construções em uma Instalação Industrial. Arranjo e Detalhamento de Tubulações.Sistemas Especiais de Tubulação. Suportes de Tubulação. Montagem e Teste de Tubulações.Visita Técnica Desenhos de TubulaçõesDesenhos de Tubulações - ExercíciosCálculo do diâmetro das tubulaçõesA Tubulação Considerada como Elemento Estrutural Cálculo da Espessura de Parede de Tubos e do Vão entre Suportes.Dilatação Térmica e Flexibilidade de Tubulações. Cálculo de Flexibilidade.Cálculo de Flexibilidade."
$ws.Range("C17").Value2 = $ws.Range("B17").Value2

# --- Método (row 20): description of teaching method instead of the
# stray second instructor's name.
$ws.Range("B20").Value2 = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula. discussão de castos práticos, visitas técnicas"
$ws.Range("C20").Value2 = $ws.Range("B20").Value2

# --- Critério (row 21): grading-criteria text instead of the method
# text that had leaked one row early.
$ws.Range("B21").Value2 = "Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula."
$ws.Range("C21").Value2 = $ws.Range("B21").Value2

# --- Norma de recuperação (row 22): make-up exam rule text.
$ws.Range("B22").Value2 = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
$ws.Range("C22").Value2 = $ws.Range("B22").Value2

# --- Bibliografia (row 23): actual bibliography listing instead of the
# make-up exam rule text that had leaked one row early.
$ws.Range("B23").Value2 = "1)TUBULAÇÕES INDUSTRIAIS - Volume I e IISilva Telles, Pedro c. - Ed. Livros Técnicos e Científicos Editora S/A2)TABELAS E GRÁFICOS PARA PROJETOS DE TUBULAÇÕES INDUSTRIAISSilva Telles, P.C./Paula Barros, Darcy G. - Ed. Interciência Ltda3)TUBULAÇÕESSilva, Remi Benedito - Editora Grêmio Politécnico da USP4)MATERIAIS PARA EQUIPAMENTOS DE PROCESSOSilva Telles, Pedro C. - Ed. Interciência Ltda5)CATÁLOGOS DIVERSOS"
$ws.Range("C23").Value2 = $ws.Range("B23").Value2
